$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Update the Runmode column (C) from "Y" to "N" for all rows except C3 and C14
# (row 14 is already "N"; row 3 is left untouched per the target diff)
$rows = @(2,4,5,6,7,8,9,10,11,12,13,15,16,17,18)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 3).Value = "N"
}

# Update the sheet view: scroll so column C is the left-most visible column,
# and move the active selection to C3
$ws.Activate()
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Range("C3").Select()
